$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I ("Location of Issue") was actually holding a date-time value all
# along -- rename the header to reflect what it really is ("Date of Last
# Update"). The underlying data in column I is left untouched.
$ws.Range("I1").Value = "Date of Last Update"

# Insert two fresh columns after the (renamed) "Date of Last Update" column,
# pushing the old J:L ("Management or Project Management Issue", "Project
# Management Issues", "Project Managing Processes & Tools") block down to
# L:N.
$ws.Columns.Item(10).Insert()
$ws.Columns.Item(10).Insert()

# New column J duplicates the TicketID column (the "weird bug" from the
# commit message). Copy the cell instead of re-typing "247" so the value
# keeps its original text type instead of being re-inferred as a number.
$ws.Range("J1").Value = "TicketID"
$ws.Range("A2").Copy($ws.Range("J2"))

# New column K restores a "Location of Issue" header, with no data recorded
# for the existing row.
$ws.Range("K1").Value = "Location of Issue"
$ws.Range("K2").Value = "None"

# Reset the selection back to the top-left cell.
$ws.Range("A1").Select() | Out-Null
